$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.657.53'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.071.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.42'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.83'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0786'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.79%  '

$ws.Range("E11").Value = '  +2.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.365.67'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.35'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.82'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.757'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.076.95'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.565.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.51'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0817'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.35'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("E24").Value = '  -1.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.138'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.87'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.42'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.33%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.61'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0622'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.30%  '

$ws.Range("E36").Value = '  +3.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.36'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.12%  '

$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0988'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.43%  '

$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.28'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0213'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.68%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.72%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.448.64'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.94%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.16'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.12%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.40'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.85%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.51'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.56%  '

$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.47%  '
